$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("E16:E19").HorizontalAlignment = -4108
